$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fix tag ordering
$ws.Range("B3").Value = "pictures,Chihuahua"

# Row 4: clear the Reference value (becomes an empty inline string)
$ws.Range("D4").Value = ""

# Row 5: new apple entry
$ws.Range("A5").Value = "C:\Users\Veeraraju_elluru\Desktop\Veeraraju\Personal\ATREE\test_images\apple.jpeg"
$ws.Range("B5").Value = "apple,tree,stem"
$ws.Range("C5").Value = "Malus pumila"
$ws.Range("D5").Value = "https://en.wikipedia.org/wiki/Apple"

# Row 6: new lotus entry (Species left blank, Reference literally "nan")
$ws.Range("A6").Value = "C:\Users\Veeraraju_elluru\Desktop\Veeraraju\Personal\ATREE\test_images\lotus.jpeg"
$ws.Range("B6").Value = "lotus,flower"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "nan"
